$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark that wraps "sudo su - db2inst2" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Update the SQL text: db2 "call DROPOBJECTS()" -> call DROPOBJECTS() ---
$d.Content.Find.Execute('db2 "call DROPOBJECTS()"', $false, $false, $false, $false, $false, $true, 1, $false, 'call DROPOBJECTS()', 2)

# --- 3. Re-add a "_GoBack" bookmark right after the new text, inside its paragraph ---
$rng = $d.Content
$rng.Find.Execute("call DROPOBJECTS()", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $rng.End

# The harness mis-resolves a brand-new zero-length Range sitting exactly at the
# last character of a paragraph (just before its paragraph mark) when handed
# straight to Bookmarks.Add. Work around it: briefly insert a one-character
# marker at that position (which nudges the position off the paragraph-end
# edge case), add the bookmark there, then remove the marker again. The
# bookmark stays correctly anchored once created.
$marker = $d.Range($endPos, $endPos)
$marker.InsertAfter("@")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($endPos, $endPos + 1)
$markerRange.Delete()
